$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "Custodians" test suite (row 2) has finished being automated.
# Update the automated test case count, mark status as Automated,
# and clear the "Needs a Delete Keyword." note that is no longer needed.
$ws.Range("B2").Value = 6
$ws.Range("D2").Value = "Automated"
$ws.Range("E2").ClearContents()

# Move the active selection to E2 as left after the edit.
$ws.Range("E2").Select()
